$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row 39 (shifts old row 40 -> 41, old row 41 -> 42; merge cells & shared
#     formulas below it are shifted/renumbered automatically by the engine) ---
$ws.Rows(39).Insert()

# --- Populate new row 39: "Critical Thickness/in" label + conversion formulas ---
$ws.Range("F39").Value = "Critical Thickness/in"
$ws.Range("G39").Formula = "=G38/25.4"
$ws.Range("H39:L39").Formula = "=H38/25.4"

# Style G39:L39 like the existing "fillId 36" band (row 41/old-34 N-column style), then apply
# the new 0.000 number format so the engine allocates a fresh (fillId36, numFmt164) cellXf.
$ws.Range("G41").Copy()
$ws.Range("G39:L39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G39:L39").NumberFormat = "0.000"

# Style E39 like the bold/white header band (G30, fontId 17), then drop its fill so the engine
# allocates a fresh (fontId 17-like, fillId0) cellXf.
$ws.Range("G30").Copy()
$ws.Range("E39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E39").Font.ThemeColor = 2
$ws.Range("E39").Interior.ColorIndex = -4142

# --- Set the sheet selection / scroll position to match the edited workbook ---
$ws.Range("G19").Select()

Write-Host "Done"
